# Update: ut 09. 03. 2021
# Revised Slovakia COVID daily-stats sheet:
#   - historical Ag-test figures (columns F/G) for rows 335-368 were
#     corrected/re-reported
#   - a brand new day (09.03.2021, Excel serial 44263) was appended as row 369

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revisions to already-existing rows (AgTests / AgPosit columns)
$updates = @(
    @{ Cell = "F335"; Value = 130530 },
    @{ Cell = "F336"; Value = 101640 },
    @{ Cell = "G336"; Value = 3311 },
    @{ Cell = "F337"; Value = 104291 },
    @{ Cell = "G337"; Value = 2959 },
    @{ Cell = "F338"; Value = 226132 },
    @{ Cell = "G338"; Value = 3183 },
    @{ Cell = "F339"; Value = 652966 },
    @{ Cell = "G339"; Value = 5459 },
    @{ Cell = "F340"; Value = 379329 },
    @{ Cell = "G340"; Value = 3246 },
    @{ Cell = "F341"; Value = 294936 },
    @{ Cell = "G341"; Value = 3676 },
    @{ Cell = "F342"; Value = 178686 },
    @{ Cell = "G342"; Value = 3075 },
    @{ Cell = "F343"; Value = 133553 },
    @{ Cell = "G343"; Value = 2958 },
    @{ Cell = "F344"; Value = 137764 },
    @{ Cell = "G344"; Value = 2533 },
    @{ Cell = "F345"; Value = 290434 },
    @{ Cell = "F346"; Value = 666446 },
    @{ Cell = "G346"; Value = 4768 },
    @{ Cell = "F347"; Value = 340573 },
    @{ Cell = "G347"; Value = 2890 },
    @{ Cell = "F348"; Value = 230699 },
    @{ Cell = "G348"; Value = 3218 },
    @{ Cell = "F350"; Value = 127558 },
    @{ Cell = "G350"; Value = 2965 },
    @{ Cell = "F351"; Value = 150205 },
    @{ Cell = "G351"; Value = 2824 },
    @{ Cell = "F352"; Value = 306570 },
    @{ Cell = "G352"; Value = 3552 },
    @{ Cell = "F353"; Value = 717228 },
    @{ Cell = "G353"; Value = 5250 },
    @{ Cell = "F354"; Value = 304563 },
    @{ Cell = "G354"; Value = 2783 },
    @{ Cell = "F355"; Value = 221736 },
    @{ Cell = "G355"; Value = 3441 },
    @{ Cell = "F356"; Value = 160379 },
    @{ Cell = "G356"; Value = 2893 },
    @{ Cell = "F357"; Value = 138425 },
    @{ Cell = "G357"; Value = 3022 },
    @{ Cell = "F358"; Value = 157542 },
    @{ Cell = "G358"; Value = 2598 },
    @{ Cell = "F359"; Value = 320034 },
    @{ Cell = "G359"; Value = 3347 },
    @{ Cell = "F360"; Value = 737831 },
    @{ Cell = "G360"; Value = 5032 },
    @{ Cell = "F361"; Value = 329441 },
    @{ Cell = "G361"; Value = 2582 },
    @{ Cell = "F362"; Value = 223397 },
    @{ Cell = "F363"; Value = 184942 },
    @{ Cell = "F364"; Value = 163442 },
    @{ Cell = "G364"; Value = 2406 },
    @{ Cell = "F365"; Value = 177359 },
    @{ Cell = "G365"; Value = 2353 },
    @{ Cell = "F366"; Value = 324988 },
    @{ Cell = "G366"; Value = 3171 },
    @{ Cell = "F367"; Value = 719768 },
    @{ Cell = "G367"; Value = 3650 },
    @{ Cell = "F368"; Value = 329552 },
    @{ Cell = "G368"; Value = 2179 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Append the new day's row (09.03.2021 -> Excel serial 44263)
$ws.Range("A369").Value = 44263
$ws.Range("B369").Value = 0
$ws.Range("C369").Value = -2126553
$ws.Range("D369").Value = -323786
$ws.Range("E369").Value = 8037
$ws.Range("F369").Value = 191371
$ws.Range("G369").Value = 2114
